$d = $word.ActiveDocument

# The first table in the document (Column1/Column2/Column3 ... Row1..Row6)
# currently has no borders defined (<w:tblPr/>). Add single-line borders
# (size 16 = 2pt in eighths-of-a-point, auto color) on every edge, inside
# and out, matching <w:tblBorders> with top/left/bottom/right/insideH/insideV.
$t = $d.Tables.Item(1)

$t.Borders.Enable = 1

$t.Borders.OutsideLineStyle = 1
$t.Borders.OutsideLineWidth = 8

$t.Borders.InsideLineStyle = 1
$t.Borders.InsideLineWidth = 8
